# Update the "F" column ("想去人数" / want-to-go headcount) values across the
# "展览", "演出" and "全部类型" sheets to match the refreshed scrape snapshot.
$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")      # sheet1 - Exhibition
$wsShow = $wb.Worksheets.Item("演出")      # sheet2 - Performance
$wsAll  = $wb.Worksheets.Item("全部类型")  # sheet4 - All types

# --- 展览 (sheet1) ---
$wsExpo.Range("F2").Value = 0
$wsExpo.Range("F6").Value = 146
$wsExpo.Range("F9").Value = 91
$wsExpo.Range("F10").Value = 738
$wsExpo.Range("F12").Value = 1136
$wsExpo.Range("F13").Value = 0
$wsExpo.Range("F15").Value = 0
$wsExpo.Range("F17").Value = 142
$wsExpo.Range("F19").Value = 0
$wsExpo.Range("F20").Value = 6208
$wsExpo.Range("F23").Value = 0
$wsExpo.Range("F24").Value = 534
$wsExpo.Range("F26").Value = 3949
$wsExpo.Range("F27").Value = 395
$wsExpo.Range("F28").Value = 36
$wsExpo.Range("F29").Value = 2551
$wsExpo.Range("F30").Value = 0
$wsExpo.Range("F33").Value = 0
$wsExpo.Range("F36").Value = 164
$wsExpo.Range("F37").Value = 0
$wsExpo.Range("F40").Value = 0
$wsExpo.Range("F41").Value = 56
$wsExpo.Range("F42").Value = 485
$wsExpo.Range("F43").Value = 0

# --- 演出 (sheet2) ---
$wsShow.Range("F2").Value = 0

# --- 全部类型 (sheet4) ---
$wsAll.Range("F4").Value = 0
$wsAll.Range("F5").Value = 0
$wsAll.Range("F8").Value = 107
$wsAll.Range("F9").Value = 0
$wsAll.Range("F11").Value = 0
$wsAll.Range("F12").Value = 216
$wsAll.Range("F13").Value = 1136
$wsAll.Range("F15").Value = 259
$wsAll.Range("F16").Value = 0
$wsAll.Range("F18").Value = 0
$wsAll.Range("F19").Value = 0
$wsAll.Range("F20").Value = 3881
$wsAll.Range("F21").Value = 6208
$wsAll.Range("F22").Value = 0
$wsAll.Range("F23").Value = 0
$wsAll.Range("F26").Value = 0
$wsAll.Range("F27").Value = 0
$wsAll.Range("F28").Value = 0
$wsAll.Range("F30").Value = 0
$wsAll.Range("F32").Value = 526
$wsAll.Range("F33").Value = 137
$wsAll.Range("F35").Value = 296
$wsAll.Range("F36").Value = 0
$wsAll.Range("F38").Value = 1556
$wsAll.Range("F39").Value = 944
$wsAll.Range("F40").Value = 43
$wsAll.Range("F45").Value = 73
$wsAll.Range("F46").Value = 574

